$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (14.42578125 -> 13.7109375 in the OOXML "width" units).
# The COM ColumnWidth property is expressed in characters and this host
# quantizes it to 1/6-character steps, so 12.8 lands on the closest
# reachable width to the target.
$ws.Columns.Item(1).ColumnWidth = 12.8

# Update cell values (column A)
$ws.Range("A1").Value = 0.071372064694537246
$ws.Range("A2").Value = 0.035047236094488086
$ws.Range("A3").Value = -0.037721326101305602
$ws.Range("A4").Value = -0.023392593344471208

# Update cell values (column B)
$ws.Range("B1").Value = -0.071372065184316297
$ws.Range("B2").Value = -0.035047236596154205
$ws.Range("B3").Value = 0.037721325593424941
$ws.Range("B4").Value = 0.023392592821762634
